$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.316.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "'1.858.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'0.7029"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.07883"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'24.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.38%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'1.871.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "'5.212"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'0.7078"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'89.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "'29.327.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'5.790"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.000007813"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'237.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("B21").Value = "'Dai"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("B22").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'2.115.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'161.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'0.1414"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'1.914"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").Value = "'1.392"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").Value = "'1.483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").Value = "'4.040"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.05189"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "'0.7105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'2.674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'0.01849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'2.681"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'1.141.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'0.9212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "'5.965"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").Value = "'0.4247"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "'70.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'102.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'0.5316"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").Value = "'1.748"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "'9.196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'7.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "
